$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shuffled match records (rows 22-24, 34-36, 48-49, 51-52, 119-120) ---
# Row 22
$ws.Range("B22").Value = 6014977
$ws.Range("F22").Value = "Moss"
$ws.Range("G22").Value = "Kongsvinger"
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = "A"
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 3.4
$ws.Range("M22").Value = 2.3
$ws.Range("N22").Value = 2.9
$ws.Range("O22").Value = 3.5
$ws.Range("P22").Value = 2.375
$ws.Range("Q22").Value = 0.25
$ws.Range("R22").Value = 1.8
$ws.Range("S22").Value = 2.05
$ws.Range("T22").Value = 2.75
$ws.Range("U22").Value = 1.95
$ws.Range("V22").Value = 1.9
$ws.Range("W22").Value = -1
$ws.Range("X22").Value = -1
$ws.Range("Y22").Value = 1.375
$ws.Range("Z22").Value = -1
$ws.Range("AA22").Value = 1.05
$ws.Range("AB22").Value = -1
$ws.Range("AC22").Value = 0.8999999999999999
# Row 23
$ws.Range("B23").Value = 6011584
$ws.Range("F23").Value = "Asane"
$ws.Range("G23").Value = "Kristiansund BK"
$ws.Range("H23").Value = 2
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = "D"
$ws.Range("K23").Value = 3.75
$ws.Range("L23").Value = 3.5
$ws.Range("M23").Value = 1.95
$ws.Range("N23").Value = 3.5
$ws.Range("O23").Value = 3.75
$ws.Range("P23").Value = 1.95
$ws.Range("Q23").Value = 0.5
$ws.Range("R23").Value = 1.875
$ws.Range("S23").Value = 1.975
$ws.Range("T23").Value = 3.25
$ws.Range("U23").Value = 1.975
$ws.Range("V23").Value = 1.875
$ws.Range("W23").Value = -1
$ws.Range("X23").Value = 2.75
$ws.Range("Y23").Value = -1
$ws.Range("Z23").Value = 0.875
$ws.Range("AA23").Value = -1
$ws.Range("AB23").Value = 0.9750000000000001
$ws.Range("AC23").Value = -1
# Row 24
$ws.Range("B24").Value = 6011943
$ws.Range("F24").Value = "KFUM"
$ws.Range("G24").Value = "Sandnes Ulf"
$ws.Range("H24").Value = 3
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = "H"
$ws.Range("K24").Value = 1.85
$ws.Range("L24").Value = 3.6
$ws.Range("M24").Value = 4
$ws.Range("N24").Value = 1.5
$ws.Range("O24").Value = 4.333
$ws.Range("P24").Value = 6
$ws.Range("Q24").Value = -1.25
$ws.Range("R24").Value = 2.025
$ws.Range("S24").Value = 1.825
$ws.Range("T24").Value = 3.25
$ws.Range("U24").Value = 1.825
$ws.Range("V24").Value = 2.025
$ws.Range("W24").Value = 0.5
$ws.Range("X24").Value = -1
$ws.Range("Y24").Value = -1
$ws.Range("Z24").Value = 1.025
$ws.Range("AA24").Value = -1
$ws.Range("AB24").Value = 0.825
$ws.Range("AC24").Value = -1
# Row 34
$ws.Range("B34").Value = 6011665
$ws.Range("F34").Value = "Ranheim"
$ws.Range("G34").Value = "Sandnes Ulf"
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 5
$ws.Range("J34").Value = "A"
$ws.Range("K34").Value = 2
$ws.Range("L34").Value = 3.75
$ws.Range("M34").Value = 3.4
$ws.Range("N34").Value = 2.1
$ws.Range("O34").Value = 3.75
$ws.Range("P34").Value = 3.3
$ws.Range("Q34").Value = -0.5
$ws.Range("R34").Value = 1.925
$ws.Range("S34").Value = 1.925
$ws.Range("T34").Value = 3.25
$ws.Range("U34").Value = 1.825
$ws.Range("V34").Value = 2.025
$ws.Range("W34").Value = -1
$ws.Range("X34").Value = -1
$ws.Range("Y34").Value = 2.3
$ws.Range("Z34").Value = -1
$ws.Range("AA34").Value = 0.925
$ws.Range("AB34").Value = 0.825
$ws.Range("AC34").Value = -1
# Row 35
$ws.Range("B35").Value = 6010976
$ws.Range("F35").Value = "Hodd"
$ws.Range("G35").Value = "KFUM"
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = "A"
$ws.Range("K35").Value = 3.2
$ws.Range("L35").Value = 3.3
$ws.Range("M35").Value = 2.2
$ws.Range("N35").Value = 4
$ws.Range("O35").Value = 3.5
$ws.Range("P35").Value = 1.909
$ws.Range("Q35").Value = 0.5
$ws.Range("R35").Value = 1.95
$ws.Range("S35").Value = 1.9
$ws.Range("T35").Value = 2.5
$ws.Range("U35").Value = 1.85
$ws.Range("V35").Value = 2
$ws.Range("W35").Value = -1
$ws.Range("X35").Value = -1
$ws.Range("Y35").Value = 0.909
$ws.Range("Z35").Value = -1
$ws.Range("AA35").Value = 0.8999999999999999
$ws.Range("AB35").Value = -1
$ws.Range("AC35").Value = 1
# Row 36
$ws.Range("B36").Value = 6011588
$ws.Range("F36").Value = "Sogndal"
$ws.Range("G36").Value = "Jerv"
$ws.Range("H36").Value = 5
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = "H"
$ws.Range("K36").Value = 1.8
$ws.Range("L36").Value = 3.6
$ws.Range("M36").Value = 4.2
$ws.Range("N36").Value = 1.533
$ws.Range("O36").Value = 4.5
$ws.Range("P36").Value = 5.75
$ws.Range("Q36").Value = -1
$ws.Range("R36").Value = 1.825
$ws.Range("S36").Value = 2.025
$ws.Range("T36").Value = 3.25
$ws.Range("U36").Value = 2
$ws.Range("V36").Value = 1.85
$ws.Range("W36").Value = 0.5329999999999999
$ws.Range("X36").Value = -1
$ws.Range("Y36").Value = -1
$ws.Range("Z36").Value = 0.825
$ws.Range("AA36").Value = -1
$ws.Range("AB36").Value = 1
$ws.Range("AC36").Value = -1
# Row 48
$ws.Range("B48").Value = 6011594
# Row 49
$ws.Range("B49").Value = 6011596
# Row 48
$ws.Range("F48").Value = "Kongsvinger"
# Row 49
$ws.Range("F49").Value = "Sandnes Ulf"
# Row 48
$ws.Range("G48").Value = "Raufoss"
# Row 49
$ws.Range("G49").Value = "Sogndal"
# Row 48
$ws.Range("H48").Value = 1
# Row 49
$ws.Range("H49").Value = 1
# Row 48
$ws.Range("I48").Value = 3
# Row 49
$ws.Range("I49").Value = 2
# Row 48
$ws.Range("J48").Value = "A"
# Row 49
$ws.Range("J49").Value = "A"
# Row 48
$ws.Range("K48").Value = 1.65
# Row 49
$ws.Range("K49").Value = 2.7
# Row 48
$ws.Range("L48").Value = 4
# Row 49
$ws.Range("L49").Value = 3.4
# Row 48
$ws.Range("M48").Value = 4.333
# Row 49
$ws.Range("M49").Value = 2.3
# Row 48
$ws.Range("N48").Value = 1.45
# Row 49
$ws.Range("N49").Value = 2.55
# Row 48
$ws.Range("O48").Value = 4.5
# Row 49
$ws.Range("O49").Value = 3.8
# Row 48
$ws.Range("P48").Value = 6.5
# Row 49
$ws.Range("P49").Value = 2.5
# Row 48
$ws.Range("Q48").Value = -1.25
# Row 49
$ws.Range("Q49").Value = 0
# Row 48
$ws.Range("R48").Value = 2.05
# Row 49
$ws.Range("R49").Value = 2
# Row 48
$ws.Range("S48").Value = 1.8
# Row 49
$ws.Range("S49").Value = 1.85
# Row 48
$ws.Range("T48").Value = 2.75
# Row 49
$ws.Range("T49").Value = 3.25
# Row 48
$ws.Range("U48").Value = 1.825
# Row 49
$ws.Range("U49").Value = 1.925
# Row 48
$ws.Range("V48").Value = 2.025
# Row 49
$ws.Range("V49").Value = 1.925
# Row 48
$ws.Range("W48").Value = -1
# Row 49
$ws.Range("W49").Value = -1
# Row 48
$ws.Range("X48").Value = -1
# Row 49
$ws.Range("X49").Value = -1
# Row 48
$ws.Range("Y48").Value = 5.5
# Row 49
$ws.Range("Y49").Value = 1.5
# Row 48
$ws.Range("Z48").Value = -1
# Row 49
$ws.Range("Z49").Value = -1
# Row 48
$ws.Range("AA48").Value = 0.8
# Row 49
$ws.Range("AA49").Value = 0.8500000000000001
# Row 48
$ws.Range("AB48").Value = 0.825
# Row 49
$ws.Range("AB49").Value = -0.5
# Row 48
$ws.Range("AC48").Value = -1
# Row 49
$ws.Range("AC49").Value = 0.4625
# Row 51
$ws.Range("B51").Value = 6011601
# Row 52
$ws.Range("B52").Value = 6011598
# Row 51
$ws.Range("F51").Value = "IK Start"
# Row 52
$ws.Range("F52").Value = "Ranheim"
# Row 51
$ws.Range("G51").Value = "Skeid"
# Row 52
$ws.Range("G52").Value = "Kongsvinger"
# Row 51
$ws.Range("H51").Value = 3
# Row 52
$ws.Range("H52").Value = 1
# Row 51
$ws.Range("I51").Value = 0
# Row 52
$ws.Range("I52").Value = 4
# Row 51
$ws.Range("J51").Value = "H"
# Row 52
$ws.Range("J52").Value = "A"
# Row 51
$ws.Range("K51").Value = 1.666
# Row 52
$ws.Range("K52").Value = 2.8
# Row 51
$ws.Range("L51").Value = 4
# Row 52
$ws.Range("L52").Value = 3.4
# Row 51
$ws.Range("M51").Value = 4
# Row 52
$ws.Range("M52").Value = 2.2
# Row 51
$ws.Range("N51").Value = 1.45
# Row 52
$ws.Range("N52").Value = 3.3
# Row 51
$ws.Range("O51").Value = 4.75
# Row 52
$ws.Range("O52").Value = 3.6
# Row 51
$ws.Range("P51").Value = 6.5
# Row 52
$ws.Range("P52").Value = 2.1
# Row 51
$ws.Range("Q51").Value = -1.25
# Row 52
$ws.Range("Q52").Value = 0.25
# Row 51
$ws.Range("R51").Value = 1.95
# Row 52
$ws.Range("R52").Value = 1.975
# Row 51
$ws.Range("S51").Value = 1.9
# Row 52
$ws.Range("S52").Value = 1.875
# Row 51
$ws.Range("T51").Value = 3.5
# Row 52
$ws.Range("T52").Value = 3
# Row 51
$ws.Range("U51").Value = 2
# Row 52
$ws.Range("U52").Value = 1.9
# Row 51
$ws.Range("V51").Value = 1.85
# Row 52
$ws.Range("V52").Value = 1.95
# Row 51
$ws.Range("W51").Value = 0.45
# Row 52
$ws.Range("W52").Value = -1
# Row 51
$ws.Range("X51").Value = -1
# Row 52
$ws.Range("X52").Value = -1
# Row 51
$ws.Range("Y51").Value = -1
# Row 52
$ws.Range("Y52").Value = 1.1
# Row 51
$ws.Range("Z51").Value = 0.95
# Row 52
$ws.Range("Z52").Value = -1
# Row 51
$ws.Range("AA51").Value = -1
# Row 52
$ws.Range("AA52").Value = 0.875
# Row 51
$ws.Range("AB51").Value = -1
# Row 52
$ws.Range("AB52").Value = 0.8999999999999999
# Row 51
$ws.Range("AC51").Value = 0.8500000000000001
# Row 52
$ws.Range("AC52").Value = -1
# Row 119
$ws.Range("B119").Value = 6014965
# Row 120
$ws.Range("B120").Value = 6011630
# Row 119
$ws.Range("F119").Value = "Moss"
# Row 120
$ws.Range("F120").Value = "Asane"
# Row 119
$ws.Range("G119").Value = "Hodd"
# Row 120
$ws.Range("G120").Value = "Ranheim"
# Row 119
$ws.Range("H119").Value = 3
# Row 120
$ws.Range("H120").Value = 3
# Row 119
$ws.Range("I119").Value = 0
# Row 120
$ws.Range("I120").Value = 0
# Row 119
$ws.Range("J119").Value = "H"
# Row 120
$ws.Range("J120").Value = "H"
# Row 119
$ws.Range("K119").Value = 1.95
# Row 120
$ws.Range("K120").Value = 2.05
# Row 119
$ws.Range("L119").Value = 3.6
# Row 120
$ws.Range("L120").Value = 3.75
# Row 119
$ws.Range("M119").Value = 3.4
# Row 120
$ws.Range("M120").Value = 3
# Row 119
$ws.Range("N119").Value = 2.15
# Row 120
$ws.Range("N120").Value = 1.95
# Row 119
$ws.Range("O119").Value = 3.5
# Row 120
$ws.Range("O120").Value = 4.2
# Row 119
$ws.Range("P119").Value = 3.4
# Row 120
$ws.Range("P120").Value = 3.2
# Row 119
$ws.Range("Q119").Value = -0.25
# Row 120
$ws.Range("Q120").Value = -0.5
# Row 119
$ws.Range("R119").Value = 1.825
# Row 120
$ws.Range("R120").Value = 2
# Row 119
$ws.Range("S119").Value = 2.025
# Row 120
$ws.Range("S120").Value = 1.85
# Row 119
$ws.Range("T119").Value = 2.5
# Row 120
$ws.Range("T120").Value = 3.25
# Row 119
$ws.Range("U119").Value = 2.05
# Row 120
$ws.Range("U120").Value = 2.025
# Row 119
$ws.Range("V119").Value = 1.8
# Row 120
$ws.Range("V120").Value = 1.825
# Row 119
$ws.Range("W119").Value = 1.15
# Row 120
$ws.Range("W120").Value = 0.95
# Row 119
$ws.Range("X119").Value = -1
# Row 120
$ws.Range("X120").Value = -1
# Row 119
$ws.Range("Y119").Value = -1
# Row 120
$ws.Range("Y120").Value = -1
# Row 119
$ws.Range("Z119").Value = 0.825
# Row 120
$ws.Range("Z120").Value = 1
# Row 119
$ws.Range("AA119").Value = -1
# Row 120
$ws.Range("AA120").Value = -1
# Row 119
$ws.Range("AB119").Value = 1.05
# Row 120
$ws.Range("AB120").Value = -0.5
# Row 119
$ws.Range("AC119").Value = -1
# Row 120
$ws.Range("AC120").Value = 0.4125

# --- Append new match rows 177-179 ---
# Row 177
$ws.Range("A177").Value = 175
$ws.Range("B177").Value = 7628522
$ws.Range("C177").Value = "Norway Division 1"
$ws.Range("D177").Value = "Norway Division 1"
$ws.Range("E177").Value = 45397.58333333334
$ws.Range("F177").Value = "Moss"
$ws.Range("G177").Value = "Ranheim"
$ws.Range("K177").Value = 2.375
$ws.Range("L177").Value = 3.75
$ws.Range("M177").Value = 2.55
$ws.Range("N177").Value = 2.1
$ws.Range("O177").Value = 3.8
$ws.Range("P177").Value = 2.875
$ws.Range("Q177").Value = -0.25
$ws.Range("R177").Value = 1.925
$ws.Range("S177").Value = 1.925
$ws.Range("T177").Value = 2.75
$ws.Range("U177").Value = 1.825
$ws.Range("V177").Value = 2.025
$ws.Range("W177").Value = 0
$ws.Range("X177").Value = 0
$ws.Range("Y177").Value = 0
$ws.Range("Z177").Value = 0
$ws.Range("AA177").Value = 0
# Row 178
$ws.Range("A178").Value = 176
$ws.Range("B178").Value = 7629079
$ws.Range("C178").Value = "Norway Division 1"
$ws.Range("D178").Value = "Norway Division 1"
$ws.Range("E178").Value = 45397.58333333334
$ws.Range("F178").Value = "Sandnes Ulf"
$ws.Range("G178").Value = "Raufoss"
$ws.Range("K178").Value = 2.4
$ws.Range("L178").Value = 3.4
$ws.Range("M178").Value = 2.7
$ws.Range("N178").Value = 2.4
$ws.Range("O178").Value = 3.4
$ws.Range("P178").Value = 2.7
$ws.Range("Q178").Value = 0
$ws.Range("R178").Value = 1.8
$ws.Range("S178").Value = 2.05
$ws.Range("T178").Value = 2.75
$ws.Range("U178").Value = 1.925
$ws.Range("V178").Value = 1.925
$ws.Range("W178").Value = 0
$ws.Range("X178").Value = 0
$ws.Range("Y178").Value = 0
$ws.Range("Z178").Value = 0
$ws.Range("AA178").Value = 0
# Row 179
$ws.Range("A179").Value = 177
$ws.Range("B179").Value = 7629016
$ws.Range("C179").Value = "Norway Division 1"
$ws.Range("D179").Value = "Norway Division 1"
$ws.Range("E179").Value = 45397.58333333334
$ws.Range("F179").Value = "Asane"
$ws.Range("G179").Value = "Sogndal"
$ws.Range("K179").Value = 2.2
$ws.Range("L179").Value = 3.6
$ws.Range("M179").Value = 2.8
$ws.Range("N179").Value = 2.2
$ws.Range("O179").Value = 3.6
$ws.Range("P179").Value = 2.8
$ws.Range("Q179").Value = -0.25
$ws.Range("R179").Value = 2
$ws.Range("S179").Value = 1.85
$ws.Range("T179").Value = 3
$ws.Range("U179").Value = 2.025
$ws.Range("V179").Value = 1.825
$ws.Range("W179").Value = 0
$ws.Range("X179").Value = 0
$ws.Range("Y179").Value = 0
$ws.Range("Z179").Value = 0
$ws.Range("AA179").Value = 0

# --- Copy styles for new rows (A: bold border style, E: date style) ---
$ws.Range("A176").Copy()
$ws.Range("A177:A179").PasteSpecial(-4122)
$ws.Range("E176").Copy()
$ws.Range("E177:E179").PasteSpecial(-4122)
